$d = $word.ActiveDocument

# Template range carrying "bold + bold-complex-script" run formatting (b + bCs),
# used to clone properly-formatted bold runs via FormattedText (Font.Bold alone
# does not reliably stamp bCs on partial-range new runs in this host).
$boldTemplate = $d.Paragraphs.Item(18).Range.Duplicate
$null = $boldTemplate.Find.Execute("Base de Datos I")

function New-BoldRun([string]$text) {
    $ft = $boldTemplate.FormattedText
    $ft.Text = $text
    return $ft
}

# --- Change 1: bullet item "Manejo de permisos a nivel de usuarios"
#     -> "Manejo de transacciones y transacciones anidadas" ---
$d.Paragraphs.Item(20).Range.Text = "Manejo de transacciones y transacciones anidadas"

# --- Change 2: "Procedimientos y funciones almacenadas" -> add trailing "." ---
$d.Paragraphs.Item(21).Range.Text = "Procedimientos y funciones almacenadas."

# --- Change 3: "Optimización de consultas a través de índices" -> add trailing "." ---
$d.Paragraphs.Item(22).Range.Text = "Optimización de consultas a través de índices."

# --- Change 4: "Un tema adicional designado por el equipo docente (cambiar cuando
#     tengamos el tema)" -> "Manejo de permisos a nivel de usuarios de base de datos." ---
$p4 = $d.Paragraphs.Item(23)
$r4 = $p4.Range
$d.Range($r4.Start, $r4.End - 1).Text = "Manejo de permisos a nivel de usuarios de base de datos."

# --- Change 5: "La gestión de permisos de usuarios es fundamental para resguardar la
#     integridad de la información y establecer distintos niveles de acceso según roles."
#     -> "El manejo de transacciones y transacciones anidadas, que garantiza la
#     atomicidad y coherencia de las operaciones múltiples, evitando inconsistencias
#     en caso de fallos, y permitiendo un control más granular en procesos complejos." ---
$p5 = $d.Paragraphs.Item(27)

$rng5a = $p5.Range.Duplicate
$null = $rng5a.Find.Execute("La ", $true, $false, $false, $false, $false, $true, 1, $false, "El ", 2)

$rng5b = $p5.Range.Duplicate
$null = $rng5b.Find.Execute("gestión de permisos de usuarios", $true, $false, $false, $false, $false, $true, 1, $false, "manejo de transacciones y transacciones anidadas", 2)

$rng5c = $p5.Range.Duplicate
$null = $rng5c.Find.Execute(" es fundamental para resguardar la integridad de la información y establecer distintos niveles de acceso según roles.", $true, $false, $false, $false, $false, $true, 1, $false, ", que garantiza la atomicidad y coherencia de las operaciones múltiples, evitando inconsistencias en caso de fallos, y permitiendo un control más granular en procesos complejos.", 2)

# --- Change 6: "Agregar el tema que elijamos " -> "El " + bold "manejo de permisos a
#     nivel de usuarios de base de datos" + " resulta esencial para definir roles,
#     restringir accesos y asegurar la confidencialidad de la información." + " " ---
$p6 = $d.Paragraphs.Item(30)
$p6.Range.Text = "El manejo de permisos a nivel de usuarios de base de datos resulta esencial para definir roles, restringir accesos y asegurar la confidencialidad de la información."

# NOTE: the formatted-text clone must be built *before* locating the destination
# range with Find - building it after the Find causes the assignment to append
# instead of replacing in place.
$boldPhrase6 = New-BoldRun("manejo de permisos a nivel de usuarios de base de datos")
$rng6 = $d.Paragraphs.Item(30).Range.Duplicate
$null = $rng6.Find.Execute("manejo de permisos a nivel de usuarios de base de datos")
$rng6.FormattedText = $boldPhrase6

# Trailing " " as its own run, matching the target run layout.
$p6end = $d.Paragraphs.Item(30).Range.End - 1
$d.Range($p6end, $p6end).InsertAfter(" ")
